$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has three task rows (2-4) followed by an empty row (5):
#   2: Hardware Testing - Voice and Sound IO
#   3: Set up Github repository
#   4: Link GTTS (Google Text to Speech) to System
#   5: (empty)
#
# We need to insert a new "Weekly Group Meeting" task as the new row 2,
# pushing the existing three rows down into rows 3-5 (filling what used to
# be the blank row 5), and rename the GTTS task while we're at it.

# Shift existing rows down one slot (bottom-up so we don't clobber data),
# copying the whole row (values + formatting) so number formats/styles
# (e.g. the date formatting on columns D/F) travel with the data.
$ws.Range("A4:G4").Copy($ws.Range("A5:G5"))
$ws.Range("A3:G3").Copy($ws.Range("A4:G4"))
$ws.Range("A2:G2").Copy($ws.Range("A3:G3"))

# Rename the GTTS task, now living in row 5.
$ws.Range("A5").Value = "Figure out GTTS (Google Text to Speech)"

# Populate the new row 2 with the weekly group meeting entry.
$ws.Range("A2").Value = "Weekly Group Meeting - Discuss development plan"
$ws.Range("B2").Value = "Group"
$ws.Range("C2").Value = "30 minutes"
$ws.Range("D2").Value = 44965
$ws.Range("E2").Value = 100
$ws.Range("F2").Value = 44965
$ws.Range("G2").Value = "30 minutes"
